$wb = $excel.ActiveWorkbook

# --- "Translations - Lab" sheet (xl/worksheets/sheet6.xml) ---
$lab = $wb.Worksheets.Item("Translations - Lab")

# The combined "Nový liquid ({{data.price}} puffíků)" label was split into a
# plain title plus separate wizard/simple-form fields, so the old cell just
# becomes the short title text.
$lab.Cells.Item(28, 3).Value = "Nový liquid"

# Copy row 37's formatting down into the five new rows before filling them in,
# so the new cells pick up the same (wrap-text) cell style used by the rest
# of the table instead of the workbook default.
$lab.Range("A37:C37").Copy()
$lab.Range("A38:C42").PasteSpecial(-4122)

$lab.Cells.Item(38, 1).Value = "cs"
$lab.Cells.Item(38, 2).Value = "lab.liquid.create.wizard.tab"
$lab.Cells.Item(38, 3).Value = "Kouzelník"

$lab.Cells.Item(39, 1).Value = "cs"
$lab.Cells.Item(39, 2).Value = "lab.liquid.create.common.tab"
$lab.Cells.Item(39, 3).Value = "Ruční zadání"

$lab.Cells.Item(40, 1).Value = "cs"
$lab.Cells.Item(40, 2).Value = "lab.liquid.create.simple.tab"
$lab.Cells.Item(40, 3).Value = "Zrychlené zadání"

$lab.Cells.Item(41, 1).Value = "cs"
$lab.Cells.Item(41, 2).Value = "lab.liquid.aromaId.label"
$lab.Cells.Item(41, 3).Value = "Aroma"

$lab.Cells.Item(42, 1).Value = "cs"
$lab.Cells.Item(42, 2).Value = "lab.liquid.aromaId.label.required"
$lab.Cells.Item(42, 3).Value = "Vyberte prosím aroma (příchuť) liquidu."

# --- View-state: "Translations - Common" loses focus, "Translations - Lab" gets it ---
$common = $wb.Worksheets.Item("Translations - Common")
$common.Activate()
$common.Range("C28").Select()

$lab.Activate()
$lab.Range("B41").Select()
